# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.143.09"
$ws.Range("E2").Value = "  +0.78%  "
# Row 3
$ws.Range("D3").Value = "3.507.78"
$ws.Range("E3").Value = "  +0.36%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.86"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.20"
# Row 7
$ws.Range("E7").Value = "  -0.01%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  -1.66%  "
# Row 9
$ws.Range("E9").Value = "  +3.12%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("E10").Value = "  -2.21%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("E11").Value = "  -0.27%  "
# Row 12
$ws.Range("D12").Value = "4.101.85"
$ws.Range("E12").Value = "  +0.05%  "
# Row 13
$ws.Range("E13").Value = "  +0.14%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.72"
$ws.Range("E14").Value = "  +5.60%  "
# Row 15
$ws.Range("D15").Value = "67.104.77"
$ws.Range("E15").Value = "  +0.68%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  +0.43%  "
# Row 17
$ws.Range("D17").Value = "3.501.39"
$ws.Range("E17").Value = "  +0.56%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -0.32%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.41"
$ws.Range("E19").Value = "  +2.58%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.02"
$ws.Range("E20").Value = "  -0.35%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.94"
$ws.Range("E21").Value = "  +0.40%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.43"
$ws.Range("E22").Value = "  +0.59%  "
# Row 23
$ws.Range("E23").Value = "  +0.15%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.536"
$ws.Range("E24").Value = "  +0.39%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.69"
$ws.Range("E25").Value = "  -0.49%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  +0.24%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +0.16%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("E28").Value = "  -0.17%  "
# Row 29
$ws.Range("E29").Value = "  -0.46%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.13"
$ws.Range("E30").Value = "  -3.67%  "
# Row 31
$ws.Range("E31").Value = "  -1.72%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  +0.18%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("E33").Value = "  -0.23%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.36"
$ws.Range("E34").Value = "  +0.29%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +0.76%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.54"
$ws.Range("E36").Value = "  +0.68%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("E37").Value = "  -1.91%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  -0.36%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("E39").Value = "  +0.94%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.65"
$ws.Range("E40").Value = "  +0.09%  "
# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.16"
$ws.Range("E41").Value = "  +1.44%  "
# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.817.01"
$ws.Range("E42").Value = "  +1.62%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0731"
$ws.Range("E43").Value = "  -1.17%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.12"
$ws.Range("E44").Value = "  -1.36%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  +1.14%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.45"
$ws.Range("E46").Value = "  -1.07%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0303"
$ws.Range("E47").Value = "  -2.38%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "341.22"
$ws.Range("E48").Value = "  -0.47%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -0.78%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.94"
$ws.Range("E50").Value = "  -0.06%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.44"
$ws.Range("E51").Value = "  -0.72%  "
